# Update the coin ranking snapshot (cryptos list) with refreshed price/volume data.
# Rows 12-14 also have their coins rotated: WrappedEther -> row12, Polkadot -> row13,
# WrappedliquidstakedEther2.0 -> row14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.892.44"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "1.638.81"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -1.81%  "

$ws.Range("D5").Value = "'214.78"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("D6").Value = "'0.5035"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("D8").Value = "'0.2568"
$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("D9").Value = "'0.06370"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("D10").Value = "'19.41"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").Value = "'0.07773"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.686.52"
$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.254"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.863.84"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "'0.5423"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").Value = "0.0₅7881"
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").Value = "'64.46"
$ws.Range("E17").Value = "  +1.27%  "

$ws.Range("D18").Value = "25.921.97"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'197.16"
$ws.Range("E20").Value = "  -3.58%  "

$ws.Range("D21").Value = "'4.363"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "'9.878"
$ws.Range("E22").Value = "  -1.53%  "

$ws.Range("D23").Value = "'5.954"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("D25").Value = "'1.898"
$ws.Range("E25").Value = "  -4.59%  "

$ws.Range("D26").Value = "'140.67"
$ws.Range("E26").Value = "  -1.41%  "

$ws.Range("D27").Value = "'0.1129"
$ws.Range("E27").Value = "  -2.46%  "

$ws.Range("D28").Value = "'6.810"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").Value = "'15.65"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").Value = "'1.239"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").Value = "'0.04927"
$ws.Range("E31").Value = "  -1.75%  "

$ws.Range("D32").Value = "'3.246"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").Value = "'3.187"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("D34").Value = "'1.532"
$ws.Range("E34").Value = "  -0.80%  "

$ws.Range("D35").Value = "'2.373"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").Value = "'0.8896"
$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("D37").Value = "'2.611"
$ws.Range("E37").Value = "  -3.09%  "

$ws.Range("D38").Value = "1.142.26"
$ws.Range("E38").Value = "  +1.62%  "

$ws.Range("D39").Value = "'0.5519"
$ws.Range("E39").Value = "  -2.62%  "

$ws.Range("D40").Value = "'0.01569"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  -1.91%  "

$ws.Range("D42").Value = "'5.685"
$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").Value = "'0.8110"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("D44").Value = "'99.53"
$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("E45").Value = "  +4.56%  "

$ws.Range("D46").Value = "1.774.58"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").Value = "'0.4522"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  -1.52%  "

$ws.Range("D49").Value = "'54.63"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").Value = "'0.05061"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("D51").Value = "'1.005"
$ws.Range("E51").Value = "  -1.34%  "
